$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# Row 3 across all sheets corresponds to the a99d8ed4... file, which has
# moved from "Ready for handoff" to handed-back status.
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

$zhcn.Range("B3").Value = $handedBack
$zhcn.Range("G2").Value = "2016-02-18 04:18:05"
$zhcn.Range("G3").Value = "2016-02-18 04:18:05"

$dede.Range("B3").Value = $handedBack
$dede.Range("G2").Value = "2016-02-18 04:18:29"
$dede.Range("G3").Value = "2016-02-18 04:18:29"
